$wb = $excel.ActiveWorkbook

# --- lipidomics_assay: update instrument_type description (new enum values) ---
$wsAssay = $wb.Worksheets.Item("lipidomics_assay")
$wsAssay.Range("E2").Value = "description: Type of mass spec the instrument used for lipidome measurements, enum: [LC-MS, MS/MS, Multidimensional MS, Ion Mobility MS, MALDI MS, GC-MS, High Mass Accuracy MS]"

# --- lipidomics_file: insert new "data_type" column (M) ---
$wsFile = $wb.Worksheets.Item("lipidomics_file")
$wsFile.Columns.Item(13).Insert()
$wsFile.Range("M1").Value = "data_type"
$wsFile.Range("M2").Value = "description: The type of data that this mapping file is associated with, enum: [LC-MS, MS/MS, Multidimensional MS, Ion Mobility MS, MALDI MS, GC-MS, High Mass Accuracy MS]"

# --- lipidomics_mapping_file: insert new "data_type" column (I) ---
$wsMap = $wb.Worksheets.Item("lipidomics_mapping_file")
$wsMap.Columns.Item(9).Insert()
$wsMap.Range("I1").Value = "data_type"
$wsMap.Range("I2").Value = "description: The type of data that this mapping file is associated with, enum: [LC-MS, MS/MS, Multidimensional MS, Ion Mobility MS, MALDI MS, GC-MS, High Mass Accuracy MS]"

# --- Restore view/selection state to match the saved workbook session ---
$wsSubject = $wb.Worksheets.Item("subject")
$wsSubject.Activate()
$wsSubject.Range("B26").Select()

$wsSample = $wb.Worksheets.Item("sample")
$wsSample.Activate()
$wsSample.Range("B50").Select()

$wsAssay.Activate()
$wsAssay.Range("E2").Select()

$wsFile.Activate()
$wsFile.Range("O16").Select()

$wsMap.Activate()
$wsMap.Range("G38").Select()
